# jobsitePermission.xlsx edit script
# Commit message: "Add new testcases and pages"
#
# The workbook has a single sheet ("Sheet1") with a small table in A1:C9.
# Column B, rows 2-9 each hold a large JSON blob describing a user's
# permission set. This edit adds a new "LoginVerifyType":"0" field to the
# UserInfo object of every one of those 8 JSON blobs (inserted right after
# the existing "IsAdvisor":false field), and moves the sheet's active
# selection to C24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldFragment = '"IsAdvisor":false,"LocationIds"'
$newFragment = '"IsAdvisor":false,"LoginVerifyType":"0","LocationIds"'

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = $cell.Value()
    if ($current.Contains($oldFragment)) {
        $updated = $current.Replace($oldFragment, $newFragment)
        $cell.Value = $updated
    }
}

# Update the saved selection on the sheet.
$ws.Range("C24").Select()
